$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.4
$ws.Range("P2").Value = 0.2
$ws.Range("S2").Value = 0.1
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.5
$ws.Range("J6").Value = 0.2857142857142857
$ws.Range("O6").Value = 0.07142857142857142
$ws.Range("Q6").Value = 0.07142857142857142
$ws.Range("R6").Value = 0.2142857142857143
$ws.Range("S6").Value = 0.3571428571428572
$ws.Range("B7").Value = 0.09090909090909091
$ws.Range("O7").Value = 0.09090909090909091
$ws.Range("R7").Value = 0.2727272727272727
$ws.Range("S7").Value = 0.5454545454545454
$ws.Range("B8").Value = 0.01886792452830189
$ws.Range("D8").Value = 0.03773584905660377
$ws.Range("F8").Value = 0.03773584905660377
$ws.Range("J8").Value = 0.03773584905660377
$ws.Range("Q8").Value = 0.1320754716981132
$ws.Range("R8").Value = 0.03773584905660377
$ws.Range("S8").Value = 0.6981132075471698
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("Q9").Value = 0.1818181818181818
$ws.Range("S9").Value = 0.7272727272727273
$ws.Range("B10").Value = 0.0847457627118644
$ws.Range("D10").Value = 0.01694915254237288
$ws.Range("F10").Value = 0.1016949152542373
$ws.Range("J10").Value = 0.0847457627118644
$ws.Range("O10").Value = 0.01694915254237288
$ws.Range("Q10").Value = 0.1355932203389831
$ws.Range("R10").Value = 0.0847457627118644
$ws.Range("S10").Value = 0.4745762711864407
$ws.Range("G11").Value = 0.06666666666666667
$ws.Range("J11").Value = 0.1333333333333333
$ws.Range("K11").Value = 0.2666666666666667
$ws.Range("L11").Value = 0.4
$ws.Range("S11").Value = 0.1333333333333333
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("L12").Value = 0.1428571428571428
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("F15").Value = 0.0625
$ws.Range("H15").Value = 0.25
$ws.Range("J15").Value = 0.375
$ws.Range("S15").Value = 0.3125
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("J16").Value = 0.1666666666666667
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("H17").Value = 0.1666666666666667
$ws.Range("I17").Value = 0.2222222222222222
$ws.Range("J17").Value = 0.4444444444444444
$ws.Range("K17").Value = 0.05555555555555555
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("H18").Value = 0.1538461538461539
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.2307692307692308
$ws.Range("K18").Value = 0.07692307692307693
$ws.Range("O18").Value = 0.1538461538461539
$ws.Range("S18").Value = 0.3076923076923077
$ws.Range("F19").Value = 0.008620689655172414
$ws.Range("H19").Value = 0.3793103448275862
$ws.Range("I19").Value = 0.05172413793103448
$ws.Range("J19").Value = 0.2327586206896552
$ws.Range("K19").Value = 0.06896551724137931
$ws.Range("M19").Value = 0.0603448275862069
$ws.Range("O19").Value = 0.06896551724137931
$ws.Range("S19").Value = 0.1293103448275862
